$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.58"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.74"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.52%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.942"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.24%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08297"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.786"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.34%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.503"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.87%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.961"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.69%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.15%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9311"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.62%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.23%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1948"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.42%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09546"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.47%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "7.70%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1064"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.97%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001311"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.82%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005946"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.00%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.531"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.86%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.100"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "9.03%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.55%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2603"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.92%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04422"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.71%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001260"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.84%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004422"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.04%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001192"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.81%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.15%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02817"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.56%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.13%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007927"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.32%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.47%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009092"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.66%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002103"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.45%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008780"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-23.50%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007302"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.38%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.06%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003670"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.96%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.06%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.06%"
